$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold price text that looks numeric (e.g. "211.95").
# The source workbook stores these as literal text, so force text
# formatting before assigning to stop Excel from auto-converting them
# to numbers (which would also lose the exact decimal representation).
$textCells = @("D5","D8","D9","D10","D16","D19","D23","D25","D29","D33","D42","D43","D44","D46","D48","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.653.88"
$ws.Range("E2").Value = "  -1.41%  "

$ws.Range("D3").Value = "1.595.86"
$ws.Range("E3").Value = "  -1.68%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "211.95"
$ws.Range("E5").Value = "  -1.43%  "

$ws.Range("E6").Value = "  -1.00%  "

$ws.Range("B8").Value = "Dogecoin"
$ws.Range("C8").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D8").Value = "0.0618"
$ws.Range("E8").Value = "  -1.75%  "

$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").Value = "0.247"
$ws.Range("E9").Value = "  -1.84%  "

$ws.Range("D10").Value = "19.69"
$ws.Range("E10").Value = "  -2.38%  "

$ws.Range("E11").Value = "  -1.43%  "

$ws.Range("D12").Value = "1.818.21"
$ws.Range("E12").Value = "  -1.74%  "

$ws.Range("D13").Value = "1.592.81"
$ws.Range("E13").Value = "  -1.78%  "

$ws.Range("E14").Value = "  -2.97%  "

$ws.Range("E15").Value = "  -2.98%  "

$ws.Range("D16").Value = "65.14"
$ws.Range("E16").Value = "  +0.42%  "

$ws.Range("D17").Value = "26.630.78"
$ws.Range("E17").Value = "  -1.39%  "

$ws.Range("D18").Value = "0.0₃0731"
$ws.Range("E18").Value = "  -1.84%  "

$ws.Range("D19").Value = "210.20"
$ws.Range("E19").Value = "  -1.93%  "

$ws.Range("E20").Value = "  +0.13%  "

$ws.Range("E21").Value = "  -2.40%  "

$ws.Range("E22").Value = "  -2.64%  "

$ws.Range("D23").Value = "2.31"
$ws.Range("E23").Value = "  -3.48%  "

$ws.Range("E24").Value = "  -2.13%  "

$ws.Range("D25").Value = "146.43"
$ws.Range("E25").Value = "  -1.20%  "

$ws.Range("E26").Value = "  +0.10%  "

$ws.Range("E27").Value = "  -3.09%  "

$ws.Range("E28").Value = "  -1.49%  "

$ws.Range("D29").Value = "15.35"
$ws.Range("E29").Value = "  -1.51%  "

$ws.Range("E30").Value = "  -1.56%  "

$ws.Range("E31").Value = "  -1.60%  "

$ws.Range("E32").Value = "  -3.62%  "

$ws.Range("D33").Value = "0.672"
$ws.Range("E33").Value = "  -10.23%  "

$ws.Range("E34").Value = "  -3.09%  "

$ws.Range("D35").Value = "1.295.26"
$ws.Range("E35").Value = "  -4.08%  "

$ws.Range("E36").Value = "  -0.27%  "

$ws.Range("E37").Value = "  -5.42%  "

$ws.Range("E38").Value = "  -3.00%  "

$ws.Range("E39").Value = "  -1.45%  "

$ws.Range("E41").Value = "  -0.96%  "

$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "5.39"
$ws.Range("E42").Value = "  +0.67%  "

$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D43").Value = "2.20"
$ws.Range("E43").Value = "  -1.51%  "

$ws.Range("D44").Value = "63.79"
$ws.Range("E44").Value = "  -2.29%  "

$ws.Range("D45").Value = "1.731.42"
$ws.Range("E45").Value = "  -1.72%  "

$ws.Range("D46").Value = "89.80"
$ws.Range("E46").Value = "  -0.34%  "

$ws.Range("E47").Value = "  -2.14%  "

$ws.Range("D48").Value = "0.839"
$ws.Range("E48").Value = "  -4.92%  "

$ws.Range("E49").Value = "  -2.68%  "

$ws.Range("E50").Value = "  -2.19%  "

$ws.Range("D51").Value = "7.52"
$ws.Range("E51").Value = "  -2.01%  "
